# Update gh-pages to output generated at 456a3b4
# Applies the F-column ("想去人数" / want-to-go count) refreshes and appends
# the new "广州·622排球少年only" event row to the 展览 (sheet 1) and
# 全部类型 (sheet 4) worksheets, plus the 本地生活 (sheet 3) F2 refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - F column updates
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    2  = 605
    3  = 274
    5  = 755
    6  = 408
    8  = 181
    10 = 249
    11 = 6930
    12 = 68
    14 = 531
    16 = 555
    18 = 426
    21 = 725
    23 = 185
    24 = 105
    25 = 332
    26 = 1039
    28 = 7
    29 = 1909
    30 = 536
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# New row 32 on "展览" - copy the A-column formatting (bold/bordered style)
# from row 31 onto row 32, then fill in the values. Only the single A cell is
# copied (not the whole row) so the worksheet's used range doesn't balloon
# out to column XFD. Date-looking text is written via NumberFormat "@" then
# ClearFormats so it stays literal text instead of being auto-converted to a
# serial date.
$ws1.Range("A31").Copy()
$ws1.Range("A32").PasteSpecial(-4122)

$ws1.Cells.Item(32, 1).Value = 31

$c = $ws1.Cells.Item(32, 2)
$c.NumberFormat = "@"
$c.Value = "2024-06-22"
$c.ClearFormats()

$ws1.Cells.Item(32, 3).Value = "广州·622排球少年only"
$ws1.Cells.Item(32, 4).Value = "岭南购物城内 广州OMG网红街"
$ws1.Cells.Item(32, 5).Value = "2024.06.22 10:00-06.22 17:30"
$ws1.Cells.Item(32, 6).Value = 0
$ws1.Cells.Item(32, 7).Value = 68
$ws1.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws1.Cells.Item(32, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - F2 update
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 307

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - F column updates
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    2  = 307
    3  = 605
    4  = 274
    6  = 755
    8  = 408
    10 = 181
    12 = 249
    13 = 6930
    14 = 68
    17 = 531
    19 = 555
    21 = 426
    28 = 725
    33 = 185
    34 = 105
    35 = 332
    36 = 1039
    38 = 7
    39 = 1909
    40 = 536
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

# New row 42 on "全部类型" - same event as above, appended at the end.
$ws4.Range("A41").Copy()
$ws4.Range("A42").PasteSpecial(-4122)

$ws4.Cells.Item(42, 1).Value = 41

$c4 = $ws4.Cells.Item(42, 2)
$c4.NumberFormat = "@"
$c4.Value = "2024-06-22"
$c4.ClearFormats()

$ws4.Cells.Item(42, 3).Value = "广州·622排球少年only"
$ws4.Cells.Item(42, 4).Value = "岭南购物城内 广州OMG网红街"
$ws4.Cells.Item(42, 5).Value = "2024.06.22 10:00-06.22 17:30"
$ws4.Cells.Item(42, 6).Value = 0
$ws4.Cells.Item(42, 7).Value = 68
$ws4.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws4.Cells.Item(42, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"
